$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the duplicated company-logo picture ---
# The sheet carries two overlapping "Picture 1" shapes anchored at the same
# spot (a leftover duplicate). Remove the first (older/larger) one and
# resize the remaining logo to its corrected dimensions.
$ws.Shapes.Item(1).Delete()
$logo = $ws.Shapes.Item(1)
$logo.Width = 93.77417322834646
$logo.Height = 55.06448818897638

# --- Row heights: row 4 goes back to the sheet's default height, row 5 keeps
# its height but is now flagged explicitly as a custom height ---
$ws.Rows.Item(4).RowHeight = 12.75
$ws.Rows.Item(5).RowHeight = 14.25

# --- Move the viewport / selection further down the form ---
$ws.Range("B31").Select()
